# "working on base billings"
# - Add new columns I/J on Sheet1 (base billings breakdown) in rows 6/7/9
# - Extend the df_a waterfall (Y/Z/AA) with 2Y / 3Y / TOTAL rows (16, 20, 23)
# - Add a new Sheet2 listing grouped legal-entity codes (G2 / ZCCR / ZCDR / ...)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Sheet1: update the df_a_1Y waterfall row and extend it ---
$ws1.Range("Z11").Value = 1833

$ws1.Range("Y16").Value = "df_a_2Y"
$ws1.Range("Z16").Value = 485
$ws1.Range("AA16").Formula = "=AA11-Z16"

$ws1.Range("Y20").Value = "df_a_3Y"
$ws1.Range("Z20").Value = 102
$ws1.Range("AA20").Formula = "=AA16-Z20"

$ws1.Range("Y23").Value = "TOTAL"
$ws1.Range("Z23").Formula = "=SUM(Z5:Z21)"
$ws1.Range("AA23").Formula = "=AA20"

# --- Sheet1: new "Total" / "grouped" mini-table around I6:J9 ---
$ws1.Range("I6").Value = "Total"
$ws1.Range("J6").Value = "grouped"

$ws1.Range("I7").Formula = "=F8"
$ws1.Range("J7").Value = 1403

$ws1.Range("I9").Value = 3058
$ws1.Range("J9").Value = 581

# --- Sheet2: new sheet with grouped legal-entity codes ---
$ws2 = $wb.Worksheets.Add($null, $ws1)

$ws2.Range("E11").Value = "G2"
$ws2.Range("E12").Value = "ZCCR"
$ws2.Range("E13").Value = "ZCDR"
$ws2.Range("E14").Value = "ZCPR"
$ws2.Range("E15").Value = "ZLCR"
$ws2.Range("E16").Value = "ZLDR"
$ws2.Range("E17").Value = "ZLG2"
$ws2.Range("E18").Value = "ZRG2"
$ws2.Range("E19").Value = "ZRL2"

$ws2.Range("I14").Select()

# Restore the active sheet/selection state on Sheet1 (tab stays on Sheet1)
$ws1.Activate()
$ws1.Range("L6").Select()
